# [ANV] started on xenon gamma decays
#
# This script:
#  1. Updates the "isotopes" sheet:
#     - re-enters the D (abundance*xn) and E (final nucleus A) formulas as a
#       single fill so Excel stores them as shared formulas
#     - shuffles the F (Sn (MeV)) column values
#     - moves the selected cell to F2 and clears the explicit tab selection
#       (the new sheet becomes the active tab instead)
#  2. Adds a new worksheet "132-primary gammas" right after "isotopes" and
#     fills in the starting data for the 132-Xe primary gamma decay work.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- isotopes sheet: rebuild D2:D10 / E2:E10 as shared formulas ---------
$ws1.Range("D2:D10").Formula = "=B2*C2"
$ws1.Range("E2:E10").Formula = "=A2+1"
# writing the formula over the whole range picks up B2's number format
# (style index 1); put the cells back to the default "Normal" style so the
# D column stays unformatted, as it was before.
$ws1.Range("D2:D10").Style = "Normal"

# --- isotopes sheet: new Sn (MeV) values in column F ---------------------
$ws1.Range("F2").Value = 8.9367175999999997
$ws1.Range("F3").Value = 9.2557229999999997
$ws1.Range("F4").Value = 6.604419
$ws1.Range("F5").Value = 7.6032999999999999
$ws1.Range("F6").Value = 6.4359000000000002
$ws1.Range("F7").Value = 6.9068535000000004
$ws1.Range("F8").Value = 6.359
$ws1.Range("F9").Value = 4.0255599999999996
$ws1.Range("F10").Value = 7.2460000000000004

# --- isotopes sheet: move the selection -----------------------------------
$ws1.Range("F2").Select()

# --- add the new worksheet right after "isotopes" -------------------------
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "132-primary gammas"

# Column widths to fit the longer headers (A/B/C), same as isotopes' cols.
$ws2.Columns.Item(1).ColumnWidth = 22.83203125
$ws2.Columns.Item(2).ColumnWidth = 20.66796875
$ws2.Columns.Item(3).ColumnWidth = 20.83203125

# Column headers for the gamma table (row 2)
$ws2.Range("A2").Value = "Gamma Energy (keV)"
$ws2.Range("B2").Value = "relative abundance"
$ws2.Range("C2").Value = "initial level (keV)"

# Reference note + neutron separation energy (row 1)
$ws2.Range("A1").Value = "reference: https://iopscience.iop.org/article/10.1088/0305-4616/14/9/009 "
$ws2.Range("D1").Value = "Sn:"
$ws2.Range("E1").Formula = "=8.9367176*1000"

# Initial level data (keV) that the gamma energy is computed from
$ws2.Range("C3").Value = 3242.6
$ws2.Range("C4").Value = 3181.4
$ws2.Range("C5").Value = 2873

# Gamma energy = Sn - initial level; A4:A5 filled together as one shared
# formula, A3 entered on its own (matches how the sheet was first built up).
$ws2.Range("A3").Formula = "=`$E`$1-C3"
$ws2.Range("A4:A5").Formula = "=`$E`$1-C4"

$ws2.Range("A10").Select()

$wb.Save()
